$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D24").Value = "[논문 요약 2021-03] Explanation Consistency Training"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222531464110"

$ws.Range("D44").Value = "Neural Processing Unit (NPU)의 기술 및 시장 동향"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/97"
